# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should look like the existing bold/centered/bordered
# header row, so copy the formatting from the last existing header cell
# (AC1) onto the three new header cells before writing their labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player on this roster shared the team's 85-77-0 season record.
$ws.Range("AD2:AD47").Value = 85
$ws.Range("AE2:AE47").Value = 77
$ws.Range("AF2:AF47").Value = 0
